# Scheduled runner update: refresh cached market-board price/profit figures
# (currentAveragePrice/NQ/HQ, LevePrice NQ/HQ and LeveProfit NQ/HQ columns)
# for a batch of leve rows across several crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 87.454544
$ws.Range("I33").Value = 60.25
$ws.Range("J33").Value = 103
$ws.Range("K33").Value = 60.25
$ws.Range("L33").Value = 103
$ws.Range("M33").Value = 168.75
$ws.Range("N33").Value = -561

# Row 64
$ws.Range("H64").Value = 3914.6353
$ws.Range("I64").Value = 3841.2222
$ws.Range("J64").Value = 3984.1843
$ws.Range("K64").Value = 3841.2222
$ws.Range("L64").Value = 3984.1843
$ws.Range("M64").Value = -3593.2222
$ws.Range("N64").Value = -4480.1843

# Row 67
$ws.Range("H67").Value = 3914.6353
$ws.Range("I67").Value = 3841.2222
$ws.Range("J67").Value = 3984.1843
$ws.Range("K67").Value = 3841.2222
$ws.Range("L67").Value = 3984.1843
$ws.Range("M67").Value = -2983.2222
$ws.Range("N67").Value = -5700.1843

# Row 125
$ws.Range("H125").Value = 2506.3635
$ws.Range("I125").Value = 983.5
$ws.Range("J125").Value = 3376.5715
$ws.Range("K125").Value = 8851.5
$ws.Range("L125").Value = 30389.1435
$ws.Range("M125").Value = -6391.5
$ws.Range("N125").Value = -35309.1435

# Row 137
$ws.Range("H137").Value = 32952.688
$ws.Range("I137").Value = 1487.65
$ws.Range("J137").Value = 85394.414
$ws.Range("K137").Value = 4462.950000000001
$ws.Range("L137").Value = 256183.242
$ws.Range("M137").Value = -1912.950000000001
$ws.Range("N137").Value = -261283.242

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 2446.611
$ws.Range("I45").Value = 2009.091
$ws.Range("J45").Value = 3134.1428
$ws.Range("K45").Value = 2009.091
$ws.Range("L45").Value = 3134.1428
$ws.Range("M45").Value = -1632.091
$ws.Range("N45").Value = -3888.1428

# Row 61
$ws.Range("H61").Value = 2244.1162
$ws.Range("I61").Value = 1696.6451
$ws.Range("J61").Value = 3658.4167
$ws.Range("K61").Value = 1696.6451
$ws.Range("L61").Value = 3658.4167
$ws.Range("M61").Value = -1484.6451
$ws.Range("N61").Value = -4082.4167

# Row 74
$ws.Range("H74").Value = 1914.25
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1914.25
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 1914.25
$ws.Range("N74").Value = -3662.25
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 1914.25
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1914.25
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 9571.25
$ws.Range("N77").Value = -18307.25
$ws.Range("M77").ClearContents()

# Row 102
$ws.Range("H102").Value = 3890
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

# Row 132
$ws.Range("H132").Value = 3374.5
$ws.Range("I132").Value = 2724.7778
$ws.Range("J132").Value = 4024.2222
$ws.Range("K132").Value = 8174.3334
$ws.Range("L132").Value = 12072.6666
$ws.Range("M132").Value = -5644.3334
$ws.Range("N132").Value = -17132.6666

# Row 136
$ws.Range("H136").Value = 2244.1162
$ws.Range("I136").Value = 1696.6451
$ws.Range("J136").Value = 3658.4167
$ws.Range("K136").Value = 5089.9353
$ws.Range("L136").Value = 10975.2501
$ws.Range("M136").Value = -2539.9353
$ws.Range("N136").Value = -16075.2501

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 2370.8928
$ws.Range("I134").Value = 2323.55
$ws.Range("J134").Value = 2489.25
$ws.Range("K134").Value = 6970.650000000001
$ws.Range("L134").Value = 7467.75
$ws.Range("M134").Value = -4435.650000000001
$ws.Range("N134").Value = -12537.75

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4636.8184
$ws.Range("I31").Value = 4490.737
$ws.Range("J31").Value = 4835.0713
$ws.Range("K31").Value = 4490.737
$ws.Range("L31").Value = 4835.0713
$ws.Range("M31").Value = -4195.737
$ws.Range("N31").Value = -5425.0713

# Row 33
$ws.Range("H33").Value = 18647.334
$ws.Range("I33").Value = 1746.6
$ws.Range("J33").Value = 30719.285
$ws.Range("K33").Value = 1746.6
$ws.Range("L33").Value = 30719.285
$ws.Range("M33").Value = -1367.6
$ws.Range("N33").Value = -31477.285

# Row 34
$ws.Range("H34").Value = 4636.8184
$ws.Range("I34").Value = 4490.737
$ws.Range("J34").Value = 4835.0713
$ws.Range("K34").Value = 4490.737
$ws.Range("L34").Value = 4835.0713
$ws.Range("M34").Value = -4288.737
$ws.Range("N34").Value = -5239.0713

# Row 132
$ws.Range("H132").Value = 2360.4736
$ws.Range("I132").Value = 1075.6
$ws.Range("J132").Value = 3788.111
$ws.Range("K132").Value = 3226.8
$ws.Range("L132").Value = 11364.333
$ws.Range("M132").Value = -696.7999999999997
$ws.Range("N132").Value = -16424.333

$ws = $wb.Worksheets.Item("CUL")
# Row 35
$ws.Range("H35").Value = 756.8
$ws.Range("I35").Value = 600
$ws.Range("J35").Value = 768
$ws.Range("K35").Value = 1800
$ws.Range("L35").Value = 2304
$ws.Range("M35").Value = -1512
$ws.Range("N35").Value = -2880

# Row 49
$ws.Range("H49").Value = 3000
$ws.Range("J49").Value = 3000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9312

# Row 57
$ws.Range("H57").Value = 3950
$ws.Range("I57").Value = 800
$ws.Range("K57").Value = 2400
$ws.Range("M57").Value = -1841

# Row 74
$ws.Range("H74").Value = 3000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 3000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 94
$ws.Range("H94").Value = 4372.5
$ws.Range("J94").Value = 4372.5
$ws.Range("L94").Value = 13117.5
$ws.Range("N94").Value = -14469.5

# Row 98
$ws.Range("H98").Value = 320.33334
$ws.Range("I98").Value = 330.5
$ws.Range("J98").Value = 300
$ws.Range("K98").Value = 991.5
$ws.Range("L98").Value = 900
$ws.Range("M98").Value = 506.5
$ws.Range("N98").Value = -3896

# Row 99
$ws.Range("H99").Value = 2550
$ws.Range("I99").Value = 1462.5
$ws.Range("K99").Value = 4387.5
$ws.Range("M99").Value = -2141.5

# Row 102
$ws.Range("H102").Value = 7475.5713
$ws.Range("J102").Value = 7475.5713
$ws.Range("L102").Value = 22426.7139
$ws.Range("N102").Value = -27294.7139

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 1723.7391
$ws.Range("I102").Value = 1711.1818
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1711.1818
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -89.18180000000007
$ws.Range("N102").Value = -5244

# Row 107
$ws.Range("H107").Value = 341.06668
$ws.Range("I107").Value = 264.3
$ws.Range("K107").Value = 264.3
$ws.Range("M107").Value = 1655.7

# Row 132
$ws.Range("H132").Value = 5012.591
$ws.Range("I132").Value = 5810.3335
$ws.Range("J132").Value = 4055.3
$ws.Range("K132").Value = 17431.0005
$ws.Range("L132").Value = 12165.9
$ws.Range("M132").Value = -14901.0005
$ws.Range("N132").Value = -17225.9

$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value = 900
$ws.Range("I17").Value = 900
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -730
$ws.Range("N17").ClearContents()

# Row 68
$ws.Range("H68").Value = 239557.14
$ws.Range("I68").Value = 770546.1
$ws.Range("J68").Value = 1527.5862
$ws.Range("K68").Value = 770546.1
$ws.Range("L68").Value = 1527.5862
$ws.Range("M68").Value = -769797.1
$ws.Range("N68").Value = -3025.5862

# Row 71
$ws.Range("H71").Value = 239557.14
$ws.Range("I71").Value = 770546.1
$ws.Range("J71").Value = 1527.5862
$ws.Range("K71").Value = 3852730.5
$ws.Range("L71").Value = 7637.931
$ws.Range("M71").Value = -3848986.5
$ws.Range("N71").Value = -15125.931

# Row 132
$ws.Range("H132").Value = 20135.572
$ws.Range("I132").Value = 28822.223
$ws.Range("J132").Value = 4499.6
$ws.Range("K132").Value = 86466.66900000001
$ws.Range("L132").Value = 13498.8
$ws.Range("M132").Value = -83936.66900000001
$ws.Range("N132").Value = -18558.8

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1944.0294
$ws.Range("I122").Value = 1441.2
$ws.Range("K122").Value = 4323.6
$ws.Range("M122").Value = -1873.6

